# Femacal de La Calera - Lechuga: add a new weekly block of 5 rows
# (one row per "Variedad": Conconina(o), Escarola, Francesa morada, Marina,
# Milanesa) at the top of the historical data table, pushing every existing
# record down by 5 rows. The workbook only contains one sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows right before the current first data block (row 1242),
# shifting rows 1242:1345 down to 1247:1350 and carrying their formatting
# (in particular the date-style on column D) down with them.
$ws.Range("A1242:R1246").Insert()

# Fill in the 5 new rows with the new week's data.
$ws.Range("A1242").Value = 3
$ws.Range("B1242").Value = "Femacal de La Calera"
$ws.Range("C1242").Value = "Coquimbo"
$ws.Range("D1242").Value = 44578
$ws.Range("E1242").Value = 5
$ws.Range("F1242").Value = 100112033
$ws.Range("G1242").Value = "Lechuga"
$ws.Range("H1242").Value = "Conconina(o)"
$ws.Range("I1242").Value = "Primera"
$ws.Range("J1242").Value = 150
$ws.Range("K1242").Value = 4000
$ws.Range("L1242").Value = 4500
$ws.Range("M1242").Value = 4233
$ws.Range("N1242").Value = "$/caja 10 unidades"
$ws.Range("O1242").Value = "Provincia de Quillota"
$ws.Range("P1242").Value = 423
$ws.Range("Q1242").Value = 10
$ws.Range("R1242").Value = "Hortaliza"

$ws.Range("A1243").Value = 3
$ws.Range("B1243").Value = "Femacal de La Calera"
$ws.Range("C1243").Value = "Coquimbo"
$ws.Range("D1243").Value = 44578
$ws.Range("E1243").Value = 5
$ws.Range("F1243").Value = 100112033
$ws.Range("G1243").Value = "Lechuga"
$ws.Range("H1243").Value = "Escarola"
$ws.Range("I1243").Value = "Primera"
$ws.Range("J1243").Value = 160
$ws.Range("K1243").Value = 6000
$ws.Range("L1243").Value = 6500
$ws.Range("M1243").Value = 6250
$ws.Range("N1243").Value = "$/caja 15 unidades"
$ws.Range("O1243").Value = "Provincia de Quillota"
$ws.Range("P1243").Value = 417
$ws.Range("Q1243").Value = 15
$ws.Range("R1243").Value = "Hortaliza"

$ws.Range("A1244").Value = 3
$ws.Range("B1244").Value = "Femacal de La Calera"
$ws.Range("C1244").Value = "Coquimbo"
$ws.Range("D1244").Value = 44578
$ws.Range("E1244").Value = 5
$ws.Range("F1244").Value = 100112033
$ws.Range("G1244").Value = "Lechuga"
$ws.Range("H1244").Value = "Francesa morada"
$ws.Range("I1244").Value = "Primera"
$ws.Range("J1244").Value = 70
$ws.Range("K1244").Value = 5000
$ws.Range("L1244").Value = 5000
$ws.Range("M1244").Value = 5000
$ws.Range("N1244").Value = "$/caja 18 unidades"
$ws.Range("O1244").Value = "Provincia de Quillota"
$ws.Range("P1244").Value = 278
$ws.Range("Q1244").Value = 18
$ws.Range("R1244").Value = "Hortaliza"

$ws.Range("A1245").Value = 3
$ws.Range("B1245").Value = "Femacal de La Calera"
$ws.Range("C1245").Value = "Coquimbo"
$ws.Range("D1245").Value = 44578
$ws.Range("E1245").Value = 5
$ws.Range("F1245").Value = 100112033
$ws.Range("G1245").Value = "Lechuga"
$ws.Range("H1245").Value = "Marina"
$ws.Range("I1245").Value = "Primera"
$ws.Range("J1245").Value = 70
$ws.Range("K1245").Value = 5000
$ws.Range("L1245").Value = 5000
$ws.Range("M1245").Value = 5000
$ws.Range("N1245").Value = "$/caja 18 unidades"
$ws.Range("O1245").Value = "Provincia de Quillota"
$ws.Range("P1245").Value = 278
$ws.Range("Q1245").Value = 18
$ws.Range("R1245").Value = "Hortaliza"

$ws.Range("A1246").Value = 3
$ws.Range("B1246").Value = "Femacal de La Calera"
$ws.Range("C1246").Value = "Coquimbo"
$ws.Range("D1246").Value = 44578
$ws.Range("E1246").Value = 5
$ws.Range("F1246").Value = 100112033
$ws.Range("G1246").Value = "Lechuga"
$ws.Range("H1246").Value = "Milanesa"
$ws.Range("I1246").Value = "Primera"
$ws.Range("J1246").Value = 145
$ws.Range("K1246").Value = 4000
$ws.Range("L1246").Value = 4500
$ws.Range("M1246").Value = 4259
$ws.Range("N1246").Value = "$/caja 20 unidades"
$ws.Range("O1246").Value = "Provincia de Quillota"
$ws.Range("P1246").Value = 213
$ws.Range("Q1246").Value = 20
$ws.Range("R1246").Value = "Hortaliza"
